$d = $word.ActiveDocument

# --- Helper: locate a paragraph whose text contains a given substring ---
function Find-ParagraphContaining($needle) {
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text -like "*$needle*") {
            return $para
        }
    }
    return $null
}

# ------------------------------------------------------------------
# 1) "Created By" value cell: "Apiwat Hantrakool" paragraph.
#    - add w:hint="cs" to the paragraph-mark run properties (pPr/rPr/rFonts)
#    - insert a new run "Mr. " before the existing "Apiwat" run
#    - move the "_GoBack" bookmark to sit right after "Mr. " and before "Apiwat"
# ------------------------------------------------------------------
$apiwatPara = Find-ParagraphContaining "Apiwat Hantrakool"
$apiwatRange = $apiwatPara.Range

$newParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00EA4117" w:rsidRPr="006203B0" w:rsidRDefault="007220BE" w:rsidP="00CD52C4"><w:pPr><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi" w:hint="cs"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve">Mr. </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006203B0"><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/></w:rPr><w:t>Apiwat</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="006203B0"><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="006203B0"><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/></w:rPr><w:t>Hantrakool</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$apiwatRange.InsertXML($newParaXml) | Out-Null

# ------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that used to sit after
#    "Use case ends." (it has now been superseded by the one added above).
# ------------------------------------------------------------------
$endsPara = Find-ParagraphContaining "Use case ends"
$endsRange = $endsPara.Range

$endsParaXml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="004C7648" w:rsidRPr="004C7648" w:rsidRDefault="004C7648" w:rsidP="004C7648"><w:pPr><w:pStyle w:val="a4"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:spacing w:after="160" w:line="259" w:lineRule="auto"/><w:jc w:val="both"/><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/></w:rPr></w:pPr><w:r w:rsidRPr="00A21C50"><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/></w:rPr><w:t>Use case ends</w:t></w:r><w:r w:rsidRPr="00A21C50"><w:rPr><w:rFonts w:asciiTheme="minorBidi" w:hAnsiTheme="minorBidi"/><w:sz w:val="28"/><w:cs/></w:rPr><w:t>.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@
$endsRange.InsertXML($endsParaXml) | Out-Null

Write-Host "Done"
